$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells for the team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from the
# existing last header cell (AC1) onto the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team's win/loss/tie record for every data row (2-50)
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 68   # AD - Wins
    $ws.Cells.Item($row, 31).Value = 94   # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF - Ties
}
